$p = $ppt.ActivePresentation

# The slide master currently carries the "Integral" theme colours (theme1.xml)
# and the notes master currently carries the "Office Theme" colours
# (theme2.xml). The edit swaps the two colour schemes between the two theme
# parts: the slide master's theme becomes the Office palette, and the notes
# master's theme becomes the Integral palette.

# Office theme palette (target for the slide master's theme, i.e. theme1.xml)
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Item($i).RGB = $officeColors[$i - 1]
}
